$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.062.20"
$ws.Range("E2").Value = "  -1.99%  "
$ws.Range("D3").Value = "2.435.55"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D5").Value = "'571.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.40%  "
$ws.Range("D6").Value = "'140.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "2.423.31"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "'5.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("E13").Value = "  -1.42%  "
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "61.088.82"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "2.411.19"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("D20").Value = "'7.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("D21").Value = "'324.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.98%  "
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").Value = "'6.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.32%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  -4.18%  "
$ws.Range("D26").Value = "'64.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("D27").Value = "'8.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.55%  "
$ws.Range("D28").Value = "'577.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.57%  "
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "0.0₃0915"
$ws.Range("E30").Value = "  -4.28%  "
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("E32").Value = "  -5.51%  "
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("E36").Value = "  -6.17%  "
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("D38").Value = "'150.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").Value = "'1.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.67%  "
$ws.Range("D40").Value = "'18.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("D41").Value = "'5.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.21%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "'41.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("E44").Value = "  -5.70%  "
$ws.Range("E45").Value = "  -4.82%  "
$ws.Range("D46").Value = "0.0₆0280"
$ws.Range("E46").Value = "  +22.79%  "
$ws.Range("D47").Value = "'141.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("D49").Value = "'0.595"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").Value = "'19.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("E51").Value = "  -3.68%  "
